$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 38
$ws.Range("F3").Value = 1379
$ws.Range("F4").Value = 13312
$ws.Range("F5").Value = 763
$ws.Range("F10").Value = 1910
$ws.Range("F13").Value = 20051
$ws.Range("G13").Value = 68
$ws.Range("F14").Value = 540
$ws.Range("F15").Value = 222
$ws.Range("F16").Value = 163
$ws.Range("F18").Value = 372
$ws.Range("F19").Value = 229
$ws.Range("F20").Value = 313
$ws.Range("F21").Value = 158
$ws.Range("F25").Value = 280
$ws.Range("F26").Value = 13
$ws.Range("F27").Value = 1355
$ws.Range("F28").Value = 54
$ws.Range("F29").Value = 377
$ws.Range("F30").Value = 77

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 299
$ws.Range("F4").Value = 4474
$ws.Range("F7").Value = 8
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 19

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 903
$ws.Range("F3").Value = 4422
$ws.Range("F4").Value = 94

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 903
$ws.Range("F3").Value = 38
$ws.Range("F5").Value = 1379
$ws.Range("F6").Value = 13312
$ws.Range("F7").Value = 299
$ws.Range("F8").Value = 763
$ws.Range("F9").Value = 4422
$ws.Range("F13").Value = 1910
$ws.Range("F16").Value = 94
$ws.Range("F17").Value = 20050
$ws.Range("G17").Value = 68
$ws.Range("F18").Value = 540
$ws.Range("F19").Value = 4474
$ws.Range("F20").Value = 222
$ws.Range("F23").Value = 163
$ws.Range("F26").Value = 8
$ws.Range("F30").Value = 372
$ws.Range("F32").Value = 313
$ws.Range("F33").Value = 158
$ws.Range("F40").Value = 280
$ws.Range("F41").Value = 13
$ws.Range("F42").Value = 1355
$ws.Range("F43").Value = 54
$ws.Range("F44").Value = 18
$ws.Range("F45").Value = 377
$ws.Range("F46").Value = 77
$ws.Range("F48").Value = 19
